$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.747.42"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").Value = "1.634.71"
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("E6").Value = "  -0.71%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("E9").Value = "  -0.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.93%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.638.06"
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("D14").Value = "1.859.41"
$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.553"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.48%  "

$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.75%  "

$ws.Range("D18").Value = "25.783.85"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("E20").Value = "  +1.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.68%  "

$ws.Range("E27").Value = "  -2.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0493"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("E32").Value = "  +1.20%  "

$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("E34").Value = "  +1.19%  "

$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.899"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("E37").Value = "  -1.32%  "

$ws.Range("D38").Value = "1.118.03"
$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("E39").Value = "  -1.65%  "

$ws.Range("E40").Value = "  -0.48%  "

$ws.Range("E41").Value = "  +0.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("E44").Value = "  -0.61%  "

$ws.Range("D45").Value = "1.771.04"
$ws.Range("E45").Value = "  -0.34%  "

$ws.Range("D46").Value = "0.0₆0109"
$ws.Range("E46").Value = "  -1.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.75%  "

$ws.Range("E48").Value = "  -2.18%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.09%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0502"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.28%  "
